$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(661, "parcela", "Praha-východ", "st. 110", "Nová Ves u Prahy", "stavební č. 110 Nová Ves u Prahy", 1, "číslo a KU za ním"),
    @(662, "parcela", "Praha-východ", "172/40", "Nová Ves u Prahy", "č. 172/40 Nová Ves u Prahy", 1, "číslo a KU za ním"),
    @(663, "budova", "Praha-východ", "st. 602", "Nová Ves u Prahy", "stavební č. 602 Nová Ves u Prahy (součástí je stavba č.p. 519, čst obce Nová Ves)", 1, "č. + KU"),
    @(664, "parcela", "Praha-východ", "173/122", "Nová Ves u Prahy", "č. 173/122 Nová Ves u Prahy", 1, "číslo a KU za ním"),
    @(665, "parcela", "Praha-východ", "173/127", "Nová Ves u Prahy", "č. 173/127 Nová Ves u Prahy", 1, "číslo a KU za ním"),
    @(16796, "budova", "Praha-východ", "st. 323", "Nová Ves u Prahy", "stavební č. 323 Nová Ves u Prahy (součástí je stavba č.p. 247, čst obce Nová Ves)", 1, "č. + KU"),
    @(16797, "parcela", "Praha-východ", "95/66", "Nová Ves u Prahy", "č. 95/66 Nová Ves u Prahy", 1, "číslo a KU za ním"),
    @(16798, "jednotka", "Praha-východ", "st. 298", "Nová Ves u Prahy", "jednotka č. 2650083, byt v budově č.p. 265, část obce Nová Ves, na parcele st. 298 Nová Ves u Prahy, podíl na společných částech domu a pozemku 620/1984", 1, "jednotka na parcele s prefixem st. + KU"),
    @(16799, "parcela", "Praha-východ", "95/218", "Nová Ves u Prahy", "č. 95/218 Nová Ves u Prahy", 1, "číslo a KU za ním"),
    @(16800, "parcela", "Praha-východ", "95/131", "Nová Ves u Prahy", "č. 95/131 Nová Ves u Prahy", 1, "číslo a KU za ním")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
    $row++
}
